# Scheduled market-data refresh for the Leve profit sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Updates the price/profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
# K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) for the leves
# whose market data moved since the last run. Some rows pick up a newly
# profitable HQ/NQ column (added), others lose one that's no longer
# applicable (cleared) - both are plain data, not formulas.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 17 - One for the Road
$ws.Range("H17").Value = 4493.2085
$ws.Range("J17").Value = 4493.2085
$ws.Range("L17").Value = 13479.6255
$ws.Range("N17").Value = -13815.6255

# Row 40 - Stuck in the Moment
$ws.Range("H40").Value = 2493.6
$ws.Range("I40").Value = 2407.3333
$ws.Range("J40").Value = 2623
$ws.Range("K40").Value = 2407.3333
$ws.Range("L40").Value = 2623
$ws.Range("M40").Value = -2232.3333
$ws.Range("N40").Value = -2973

# Row 43 - Growing Is Knowing
$ws.Range("H43").Value = 2880.3333
$ws.Range("I43").Value = 1392
$ws.Range("J43").Value = 3624.5
$ws.Range("K43").Value = 1392
$ws.Range("L43").Value = 3624.5
$ws.Range("M43").Value = -1323
$ws.Range("N43").Value = -3762.5

# Row 80 - Cleansing the Wicked Humours
$ws.Range("H80").Value = 8467.666999999999
$ws.Range("I80").Value = 499
$ws.Range("J80").Value = 12452
$ws.Range("K80").Value = 1497
$ws.Range("L80").Value = 37356
$ws.Range("M80").Value = -499
$ws.Range("N80").Value = -39352

# Row 83 - Washing Away the Sins (L)
$ws.Range("H83").Value = 8467.666999999999
$ws.Range("I83").Value = 499
$ws.Range("J83").Value = 12452
$ws.Range("K83").Value = 4491
$ws.Range("L83").Value = 112068
$ws.Range("M83").Value = 501
$ws.Range("N83").Value = -122052

# Row 132 - Fast-forwarding Flora
$ws.Range("H132").Value = 1019.5
$ws.Range("I132").Value = 966.2778
$ws.Range("K132").Value = 2898.8334
$ws.Range("M132").Value = -368.8334

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32 - Ingot We Trust
$ws.Range("H32").Value = 6923.619
$ws.Range("I32").Value = 6269.85
$ws.Range("K32").Value = 6269.85
$ws.Range("M32").Value = -5982.85

# Row 74 - As the Bolt Flies
$ws.Range("H74").Value = 2257.318
$ws.Range("I74").Value = 1861.9333
$ws.Range("J74").Value = 3104.5715
$ws.Range("K74").Value = 1861.9333
$ws.Range("L74").Value = 3104.5715
$ws.Range("M74").Value = -987.9332999999999
$ws.Range("N74").Value = -4852.5715

# Row 77 - Heavy Metal Banned (L)
$ws.Range("H77").Value = 2257.318
$ws.Range("I77").Value = 1861.9333
$ws.Range("J77").Value = 3104.5715
$ws.Range("K77").Value = 9309.666499999999
$ws.Range("L77").Value = 15522.8575
$ws.Range("M77").Value = -4941.666499999999
$ws.Range("N77").Value = -24258.8575

# Row 81 - A Halonic Masquerade
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

# Row 84 - Why I Wear a Mask (L)
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

# --- BSM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 20 - Smelt and Dealt
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 22 - Riveting Run
$ws.Range("H22").Value = 564.6
$ws.Range("I22").Value = 459.69232
$ws.Range("J22").Value = 759.4286
$ws.Range("K22").Value = 459.69232
$ws.Range("L22").Value = 759.4286
$ws.Range("M22").Value = -286.69232
$ws.Range("N22").Value = -1105.4286

# --- CRP ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 22 - Driving Up the Wall
$ws.Range("H22").Value = 23571.455
$ws.Range("I22").Value = 1065.3334
$ws.Range("J22").Value = 50578.8
$ws.Range("K22").Value = 1065.3334
$ws.Range("L22").Value = 50578.8
$ws.Range("M22").Value = -715.3334
$ws.Range("N22").Value = -51278.8

# Row 54 - The Turning Point
$ws.Range("H54").Value = 40748.668
$ws.Range("J54").Value = 40080
$ws.Range("L54").Value = 40080
$ws.Range("N54").Value = -41396

# Row 138 - Bow Out
$ws.Range("H138").Value = 149022.67
$ws.Range("I138").Value = 40499
$ws.Range("K138").Value = 40499
$ws.Range("M138").Value = -35359

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 14 - Keep Your Powder Dry
$ws.Range("H14").Value = 1077.4
$ws.Range("I14").Value = 1077.4
$ws.Range("K14").Value = 3232.2
$ws.Range("M14").Value = -3059.2

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 80 - Needs More Prayerbell
$ws.Range("H80").Value = 549
$ws.Range("I80").Value = 549
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 549
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 449
$ws.Range("N80").ClearContents()

# Row 83 - With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 549
$ws.Range("I83").Value = 549
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 2745
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 2247
$ws.Range("N83").ClearContents()

# Row 102 - Put the Metal to the Peddle
$ws.Range("H102").Value = 1125.6
$ws.Range("I102").Value = 1150.3
$ws.Range("J102").Value = 1076.2
$ws.Range("K102").Value = 1150.3
$ws.Range("L102").Value = 1076.2
$ws.Range("M102").Value = 471.7
$ws.Range("N102").Value = -4320.2

# Row 122 - Awarding Academic Excellence
$ws.Range("H122").Value = 6237.857
$ws.Range("I122").Value = 4931.6
$ws.Range("K122").Value = 14794.8
$ws.Range("M122").Value = -12344.8

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 16 - Saddle Sore
$ws.Range("H16").Value = 3498.6
$ws.Range("I16").Value = 3873.5
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 3873.5
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -3703.5
$ws.Range("N16").Value = -2339

# Row 20 - Choke Hold
$ws.Range("H20").Value = 199
$ws.Range("I20").Value = 199
$ws.Range("K20").Value = 199
$ws.Range("M20").Value = 27

# Row 22 - Skin off Their Backs
$ws.Range("H22").Value = 725
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 725
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 725
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1315

# Row 27 - Fire and Hide
$ws.Range("H27").Value = 725
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 725
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 725
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -939

# Row 40 - Best Served Toad
$ws.Range("H40").Value = 3673.8333
$ws.Range("I40").Value = 2959.6
$ws.Range("J40").Value = 4184
$ws.Range("K40").Value = 2959.6
$ws.Range("L40").Value = 4184
$ws.Range("M40").Value = -2823.6
$ws.Range("N40").Value = -4456

# Row 46 - Supply Side Logic
$ws.Range("H46").Value = 57311.11
$ws.Range("I46").Value = 167433.33
$ws.Range("J46").Value = 2250
$ws.Range("K46").Value = 167433.33
$ws.Range("L46").Value = 2250
$ws.Range("M46").Value = -167245.33
$ws.Range("N46").Value = -2626

# Row 82 - Trainin' the Neck
$ws.Range("H82").Value = 1890
$ws.Range("I82").Value = 1985
$ws.Range("J82").Value = 1700
$ws.Range("K82").Value = 1985
$ws.Range("L82").Value = 1700
$ws.Range("M82").Value = -1624
$ws.Range("N82").Value = -2422

# Row 85 - Training Is Only Skintight (L)
$ws.Range("H85").Value = 1890
$ws.Range("I85").Value = 1985
$ws.Range("J85").Value = 1700
$ws.Range("K85").Value = 1985
$ws.Range("L85").Value = 1700
$ws.Range("M85").Value = -737
$ws.Range("N85").Value = -4196

# Row 96 - Off the Cuff
$ws.Range("H96").Value = 49899.332
$ws.Range("J96").Value = 49899.332
$ws.Range("L96").Value = 49899.332
$ws.Range("N96").Value = -55391.332

# Row 122 - Hell on Leather
$ws.Range("H122").Value = 7410.4644
$ws.Range("I122").Value = 8111.7646
$ws.Range("K122").Value = 24335.2938
$ws.Range("M122").Value = -21885.2938

# Row 132 - Tenets of Tanning
$ws.Range("H132").Value = 3905.2222
$ws.Range("J132").Value = 4630.2
$ws.Range("L132").Value = 13890.6
$ws.Range("N132").Value = -18950.6

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 55 - A Matter of Import
$ws.Range("H55").Value = 9780.799999999999
$ws.Range("I55").Value = 4182.6665
$ws.Range("J55").Value = 18178
$ws.Range("K55").Value = 4182.6665
$ws.Range("L55").Value = 18178
$ws.Range("M55").Value = -3905.6665
$ws.Range("N55").Value = -18732

# Row 81 - Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 1251560.8
$ws.Range("I81").Value = 1495.6666
$ws.Range("J81").Value = 2001599.8
$ws.Range("K81").Value = 2991.3332
$ws.Range("L81").Value = 4003199.6
$ws.Range("M81").Value = -1930.3332
$ws.Range("N81").Value = -4005321.6

# Row 84 - To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 1251560.8
$ws.Range("I84").Value = 1495.6666
$ws.Range("J84").Value = 2001599.8
$ws.Range("K84").Value = 14956.666
$ws.Range("L84").Value = 20015998
$ws.Range("M84").Value = -9652.666000000001
$ws.Range("N84").Value = -20026606

# Row 99 - Say Yes to Formal Dress
$ws.Range("H99").Value = 190000
$ws.Range("J99").Value = 190000
$ws.Range("L99").Value = 190000
$ws.Range("N99").Value = -195990

# Row 107 - Flax Wax
$ws.Range("H107").Value = 681.7778
$ws.Range("J107").Value = 851
$ws.Range("L107").Value = 2553
$ws.Range("N107").Value = -6393

# Row 135 - In Line with Linen
$ws.Range("H135").Value = 295603.75
$ws.Range("J135").Value = 295603.75
$ws.Range("L135").Value = 295603.75
$ws.Range("N135").Value = -305743.75
